$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "Total"
$ws.Range("F16").Formula = "=SUM(F3:F12)"
$ws.Range("E17").Select() | Out-Null
